# Auto-generated edit script: adds "Date and Time" row at the top and a
# "Cycle Count of battery" row near the bottom (plus two new trailing
# speed-bucket rows), relabels several existing rows, and refreshes values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the very top; this shifts all existing rows down by
# one and carries the [hh]:mm:ss number format that was on B1 onto B2,
# matching the target file's style layout.
$ws.Rows.Item(1).Insert()

# Write every label/value pair for the final 45-row layout.
$ws.Cells.Item(1,1).Value = "Date and Time"
$ws.Cells.Item(1,2).Value = "2024-03-12 16:59:29.259000 to 2024-03-12 18:03:22.519000"
$ws.Cells.Item(2,1).Value = "Total time taken for the ride"
$ws.Cells.Item(2,2).Value = 0.04427361111111111
$ws.Cells.Item(3,1).Value = "Actual Ampere-hours (Ah)"
$ws.Cells.Item(3,2).Value = 27.28397611111111
$ws.Cells.Item(4,1).Value = "Actual Watt-hours (Wh)"
$ws.Cells.Item(4,2).Value = 1403.2693823775
$ws.Cells.Item(5,1).Value = "Starting SoC (Ah)"
$ws.Cells.Item(5,2).Value = 39.512
$ws.Cells.Item(6,1).Value = "Ending SoC (Ah)"
$ws.Cells.Item(6,2).Value = 12.605
$ws.Cells.Item(7,1).Value = "Starting SoC (%)"
$ws.Cells.Item(7,2).Value = 99.0
$ws.Cells.Item(8,1).Value = "Ending SoC (%)"
$ws.Cells.Item(8,2).Value = 31.0
$ws.Cells.Item(9,1).Value = "Total distance covered (km)"
$ws.Cells.Item(9,2).Value = 30.98848621502767
$ws.Cells.Item(10,1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10,2).Value = 45.28357315166281
$ws.Cells.Item(11,1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(11,2).Value = 68.0
$ws.Cells.Item(12,1).Value = "Mode"
$ws.Cells.Item(12,2).Value = "Custom mode`n97.62%`nEco mode`n0.11%"
$ws.Cells.Item(13,1).Value = "Peak Power(kW)"
$ws.Cells.Item(13,2).Value = 4861.082784
$ws.Cells.Item(14,1).Value = "Average Power(kW)"
$ws.Cells.Item(14,2).Value = -1326.62021443251
$ws.Cells.Item(15,1).Value = "Total Energy Regenerated(kWh)"
$ws.Cells.Item(15,2).Value = 0.01068924722222222
$ws.Cells.Item(16,1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(16,2).Value = 0.0007617329881871818
$ws.Cells.Item(17,1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(17,2).Value = 3.436
$ws.Cells.Item(18,1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(18,2).Value = 3.095
$ws.Cells.Item(19,1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19,2).Value = 0.3409999999999997
$ws.Cells.Item(20,1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20,2).Value = 38.0
$ws.Cells.Item(21,1).Value = "Maximum Temperature(C)"
$ws.Cells.Item(21,2).Value = 44.0
$ws.Cells.Item(22,1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(22,2).Value = 6.0
$ws.Cells.Item(23,1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23,2).Value = 58.0
$ws.Cells.Item(24,1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24,2).Value = 59.0
$ws.Cells.Item(25,1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25,2).Value = 56.0
$ws.Cells.Item(26,1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26,2).Value = 54.0
$ws.Cells.Item(27,1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27,2).Value = 94.0
$ws.Cells.Item(28,1).Value = "Abnormal Motor Temperature Detected(C)"
$ws.Cells.Item(28,2).Value = 0.0
$ws.Cells.Item(29,1).Value = "highest cell temp(C)"
$ws.Cells.Item(29,2).Value = 44.0
$ws.Cells.Item(30,1).Value = "lowest cell temp(C)"
$ws.Cells.Item(30,2).Value = 38.0
$ws.Cells.Item(31,1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Cells.Item(31,2).Value = 6.0
$ws.Cells.Item(32,1).Value = "Battery Voltage(V)"
$ws.Cells.Item(32,2).Value = 54.0
$ws.Cells.Item(33,1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(33,2).Value = 1.47333471
$ws.Cells.Item(34,1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(34,2).Value = 0.0000001069959847494553
$ws.Cells.Item(35,1).Value = "Cycle Count of battery"
$ws.Cells.Item(35,2).Value = 46.0
$ws.Cells.Item(36,1).Value = "Idling time percentage"
$ws.Cells.Item(36,2).Value = 2.50272034820457
$ws.Cells.Item(37,1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(37,2).Value = 24.19964492297119
$ws.Cells.Item(38,1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(38,2).Value = 7.485252849206804
$ws.Cells.Item(39,1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(39,2).Value = 11.44837065460168
$ws.Cells.Item(40,1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(40,2).Value = 14.70706145123418
$ws.Cells.Item(41,1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(41,2).Value = 20.27375293511254
$ws.Cells.Item(42,1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(42,2).Value = 16.39940438691942
$ws.Cells.Item(43,1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(43,2).Value = 2.89502319454785
$ws.Cells.Item(44,1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(44,2).Value = 0.0
$ws.Cells.Item(45,1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(45,2).Value = 0.0
